$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.7257523197931459, 3.817057569599056, 1.363035514247575, 0.7257523197931457),
    @(0.7143413178530797, 3.975964136606144, 1.462172020271649, 0.6844083110598368),
    @(0.2208274643071332, 2.572306780362553, 0.2764102973443819, 0.220827464307133),
    @(0.200350224519713, 2.640215549496584, 0.2843422165241118, 0.1754672905536122),
    @(0.4463262469701371, 0.6231076055098256, 0.4160801037814331, 0.446326246970137),
    @(0.4241068855029057, 0.6482324452444937, 0.4353507314332526, 0.3938521138762955),
    @(0.4351124807044907, 290.4257300449753, 75.9781313613845, 0.435112480704491),
    @(0.415339100968617, 300.4003350860203, 78.77358394294578, 0.3927801320263483),
    @(0.620685364174401, 5.794343485914299, 1.937925912850309, 0.6206853641744012),
    @(0.6058007273239474, 6.027106593323159, 2.030112941361703, 0.5837390607199583)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i][0]
    $ws.Cells.Item($row, 4).Value = $values[$i][1]
    $ws.Cells.Item($row, 5).Value = $values[$i][2]
    $ws.Cells.Item($row, 6).Value = $values[$i][3]
}
